# Apply the "may 9th" data update:
#  1. Insert a new row at row 2 (pushing the existing 20 data rows down by one).
#  2. Populate the new row 2 with fresh C:H sensor values (A/B continue the
#     existing timestamp/label sequence).
#  3. Append 9 brand-new rows (new timestamps 2100..2900) after the old last
#     row (now row 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above the current row 2 (shifts rows 2:21 -> 3:22) ---
$ws.Rows(2).Insert()
# Insert() copies the formatting of the row above (the bold/centered header
# style) into the freshly inserted row; the data rows in this sheet carry no
# explicit style, so strip it back off to match.
$ws.Rows(2).ClearFormats()

# --- 2. Fill the newly inserted row 2 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "struggle"
$ws.Range("C2").Value = 0.3229818344116211
$ws.Range("D2").Value = 0.6911778450012207
$ws.Range("E2").Value = 0.0410229265689849
$ws.Range("F2").Value = -0.0050396383740007
$ws.Range("G2").Value = -0.0062613687478005
$ws.Range("H2").Value = -0.0682641938328743

# --- 2b. Re-stamp the timestamp/label columns for the shifted rows (3..22)
#          so they continue the 0,100,200,... sequence (the Insert() above
#          only shifted cells down; it didn't renumber them). ---
for ($r = 3; $r -le 22; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
    $ws.Cells.Item($r, 2).Value = "struggle"
}

# --- 3. Append 9 new rows of data (rows 23..31) ---
$newRows = @(
    @(2100, 0.9749262332916433, 1.300361778587099, -6.267426431179062, -0.4335615932941437, 0.1406517177820205, -0.8185594081878662),
    @(2200, -3.844243764877326, 1.833226948976521, -1.409952521324157, 0.09498954564332961, -0.7519751191139221, -0.1093448773026466),
    @(2300, 0.4797788858413697, -0.523662269115448, -1.702465817332268, 0.1846340149641037, -1.312596678733826, 0.0687223374843597),
    @(2400, 1.155098915100098, 1.092013478279114, 1.727226853370667, 0.6478226184844971, -0.9091202020645142, -0.1838704347610473),
    @(2500, -1.098365545272828, -0.6193101108074199, 0.1845241859555233, -0.1064432710409164, -0.09178250283002851, 0.0652098655700683),
    @(2600, -0.8518145084381094, -0.03355145454406605, 0.7549576908350003, -0.042302418500185, 0.3572034537792206, 0.1937969923019409),
    @(2700, 0.405293345451355, 0.8384262472391129, 0.3231545425951481, -0.2768746614456177, 0.2338086664676666, -0.1817324161529541),
    @(2800, 0.2438197135925255, 0.4860433936119046, -0.09267929568886754, 0.0734565481543541, 0.1968513280153274, 0.1055269688367843),
    @(2900, -0.07322704792022328, 0.1344193816185026, -0.148086081258953, -0.1507309973239898, -0.0175623763352632, 0.08170322328805921)
)

$row = 23
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = "struggle"
    $ws.Cells.Item($row, 3).Value = $data[1]
    $ws.Cells.Item($row, 4).Value = $data[2]
    $ws.Cells.Item($row, 5).Value = $data[3]
    $ws.Cells.Item($row, 6).Value = $data[4]
    $ws.Cells.Item($row, 7).Value = $data[5]
    $ws.Cells.Item($row, 8).Value = $data[6]
    $row++
}
